# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.523.70'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.801.07'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.34'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.600'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.14'
$ws.Range('E8').Value = '  +13.28%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0666'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0996'
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('D12').Value = '2.061.23'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = '1.799.91'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.91'
$ws.Range('E14').Value = '  -3.03%  '
$ws.Range('D15').Value = '34.471.94'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.39'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.04'
$ws.Range('E19').Value = '  -2.39%  '
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.28'
$ws.Range('E23').Value = '  +4.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.01'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.66'
$ws.Range('E26').Value = '  -3.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.36'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.120'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('B35').Value = 'Aave'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '87.30'
$ws.Range('E35').Value = '  +7.65%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.323.43'
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.647'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0187'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.72'
$ws.Range('E40').Value = '  +11.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.33'
$ws.Range('E41').Value = '  -1.16%  '
$ws.Range('E42').Value = '  +5.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.44'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.934'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0517'
$ws.Range('E46').Value = '  +4.04%  '
$ws.Range('D47').Value = '1.963.94'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.61'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -1.21%  '
